# "fixed errors and overview now working"
#
# 1) Correct two course-name typos.
# 2) Add four new "overview" rows (25-28) with HP/course/grade/credit data,
#    styled like the rest of the table but without the inner top/bottom
#    borders (a thin divider block tacked on below the main table).
# 3) Leave the final selection on B18, matching the saved workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text corrections -------------------------------------------------
# (Edit order mirrors the original authoring session so the shared-string
# table ends up in the same append order.)
$ws.Range("B24").Value = "Introduktion till artificiella neuronnätverk och deep learning"
$ws.Range("B7").Value  = "Programmeringsteknik Fördjupningskurs"

# --- New rows: course/term labels (also in authoring order) -----------
$ws.Range("B26").Value = "Kandidatarbete"
$ws.Range("B27").Value = "Programvaruutveckling i grupp"
$ws.Range("B25").Value = "Datorer och Datoranvändning"
$ws.Range("B28").Value = "Utvärdering av Programvarusystem"

# --- New rows: remaining columns ---------------------------------------
$ws.Range("A25").Value = "HT20"
$ws.Range("C25").Value = -1
$ws.Range("D25").Value = 3

$ws.Range("A26").Value = "HT22"
$ws.Range("C26").Value = -1
$ws.Range("D26").Value = 15

$ws.Range("A27").Value = "HT21"
$ws.Range("C27").Value = -1
$ws.Range("D27").Value = 7.5

$ws.Range("A28").Value = "VT21"
$ws.Range("C28").Value = -1
$ws.Range("D28").Value = 7

# --- Formatting: reuse the existing table look for the new block ------
# Column A keeps the "medium left border" look of the rest of column A.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A25:A28").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Columns B-D keep the plain thin-box look of the rest of the table.
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B25:D28").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# The new block sits below the table without repeating the inner
# top/bottom rules, so strip those on every new cell.
foreach ($r in 25..28) {
    foreach ($col in @("A", "B", "C", "D")) {
        $cell = $ws.Range("$col$r")
        $cell.Borders.Item(8).LineStyle = -4142   # xlEdgeTop    -> none
        $cell.Borders.Item(9).LineStyle = -4142   # xlEdgeBottom -> none
    }
}

# D28 was left unformatted (default style) in the saved file.
$ws.Range("D28").ClearFormats() | Out-Null
$ws.Range("D28").Value = 7

# --- Final selection ----------------------------------------------------
$ws.Range("B18").Select() | Out-Null
